{"js": "// Post Section 3 Video 4\n// Title: \".NET 6 Entity Framework: How-To Guide for Professionals\"\n//     -> \".NET Core 6 Entity Framework: How-To Guide for Professionals\"\nconst body = context.document.body;\n\nconst titleResults = body.search(\n  \".NET 6 Entity Framework: How-To Guide for Professionals\",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\n    \".NET Core 6 Entity Framework: How-To Guide for Professionals\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// Subtitle: \" advanced and scalable architectures using .net 6 Entity Framework\"\n//       -> \" advanced and scalable architectures using .net Entity Framework Core 6.0\"\nconst subtitleResults = body.search(\n  \" advanced and scalable architectures using .net 6 Entity Framework\",\n  { matchCase: true }\n);\nsubtitleResults.load(\"items\");\nawait context.sync();\n\nif (subtitleResults.items.length > 0) {\n  subtitleResults.items[0].insertText(\n    \" advanced and scalable architectures using .net Entity Framework Core 6.0\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Post Section 3 Video 4\n$d = $word.ActiveDocument\n\n# Title: \".NET 6 Entity Framework: How-To Guide for Professionals\"\n#     -> \".NET Core 6 Entity Framework: How-To Guide for Professionals\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \".NET 6 Entity Framework: How-To Guide for Professionals\"\n$find.Replacement.Text = \".NET Core 6 Entity Framework: How-To Guide for Professionals\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Subtitle: \" advanced and scalable architectures using .net 6 Entity Framework\"\n#       -> \" advanced and scalable architectures using .net Entity Framework Core 6.0\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \" advanced and scalable architectures using .net 6 Entity Framework\"\n$find2.Replacement.Text = \" advanced and scalable architectures using .net Entity Framework Core 6.0\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
